$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 23333.334
$ws.Cells.Item(21, 9).Value = 20000
$ws.Cells.Item(21, 10).Value = 30000
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 30000
$ws.Cells.Item(21, 13).Value = -19532
$ws.Cells.Item(21, 14).Value = -30936
$ws.Cells.Item(23, 8).Value = 23333.334
$ws.Cells.Item(23, 9).Value = 20000
$ws.Cells.Item(23, 10).Value = 30000
$ws.Cells.Item(23, 11).Value = 20000
$ws.Cells.Item(23, 12).Value = 30000
$ws.Cells.Item(23, 13).Value = -19766
$ws.Cells.Item(23, 14).Value = -30468
$ws.Cells.Item(28, 8).Value = 341.8421
$ws.Cells.Item(28, 9).Value = 199.82353
$ws.Cells.Item(28, 10).Value = 1549
$ws.Cells.Item(28, 11).Value = 199.82353
$ws.Cells.Item(28, 12).Value = 1549
$ws.Cells.Item(28, 13).Value = 285.17647
$ws.Cells.Item(28, 14).Value = -2519
$ws.Cells.Item(137, 8).Value = 2846.08
$ws.Cells.Item(137, 9).Value = 2479.4375
$ws.Cells.Item(137, 11).Value = 7438.3125
$ws.Cells.Item(137, 13).Value = -4888.3125
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 2375.75
$ws.Cells.Item(13, 9).Value = 5000
$ws.Cells.Item(13, 10).Value = 1501
$ws.Cells.Item(13, 11).Value = 5000
$ws.Cells.Item(13, 12).Value = 1501
$ws.Cells.Item(13, 13).Value = -4856
$ws.Cells.Item(13, 14).Value = -1789
$ws.Cells.Item(32, 8).Value = 32471.225
$ws.Cells.Item(32, 9).Value = 31101.316
$ws.Cells.Item(32, 10).Value = 58499.5
$ws.Cells.Item(32, 11).Value = 31101.316
$ws.Cells.Item(32, 12).Value = 58499.5
$ws.Cells.Item(32, 13).Value = -30814.316
$ws.Cells.Item(32, 14).Value = -59073.5
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 13).Value = $null
$ws.Cells.Item(109, 8).Value = 50000
$ws.Cells.Item(109, 10).Value = 50000
$ws.Cells.Item(109, 12).Value = 50000
$ws.Cells.Item(109, 14).Value = -52774
$ws.Cells.Item(132, 8).Value = 8402.869000000001
$ws.Cells.Item(132, 9).Value = 6474.6665
$ws.Cells.Item(132, 11).Value = 19423.9995
$ws.Cells.Item(132, 13).Value = -16893.9995
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 107.27778
$ws.Cells.Item(7, 9).Value = 141.41667
$ws.Cells.Item(7, 10).Value = 39
$ws.Cells.Item(7, 11).Value = 141.41667
$ws.Cells.Item(7, 12).Value = 39
$ws.Cells.Item(7, 13).Value = -28.41667000000001
$ws.Cells.Item(7, 14).Value = -265
$ws.Cells.Item(105, 8).Value = 1791.1666
$ws.Cells.Item(105, 9).Value = 1149.4
$ws.Cells.Item(105, 11).Value = 1149.4
$ws.Cells.Item(105, 13).Value = 597.5999999999999
$ws.Cells.Item(122, 8).Value = 58380
$ws.Cells.Item(122, 9).Value = 78386.62
$ws.Cells.Item(122, 11).Value = 235159.86
$ws.Cells.Item(122, 13).Value = -232709.86
$ws.Cells.Item(132, 8).Value = 19143.129
$ws.Cells.Item(132, 9).Value = 741.4828
$ws.Cells.Item(132, 11).Value = 2224.4484
$ws.Cells.Item(132, 13).Value = 305.5515999999998
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 5457.615
$ws.Cells.Item(2, 9).Value = 23367.666
$ws.Cells.Item(2, 10).Value = 84.59999999999999
$ws.Cells.Item(2, 11).Value = 140205.996
$ws.Cells.Item(2, 12).Value = 507.6
$ws.Cells.Item(2, 13).Value = -140092.996
$ws.Cells.Item(2, 14).Value = -733.5999999999999
$ws.Cells.Item(5, 8).Value = 728.7143
$ws.Cells.Item(5, 9).Value = 621.4
$ws.Cells.Item(5, 11).Value = 1864.2
$ws.Cells.Item(5, 13).Value = -1752.2
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = $null
$ws.Cells.Item(17, 13).Value = $null
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(31, 8).Value = 416.66666
$ws.Cells.Item(31, 9).Value = 125
$ws.Cells.Item(31, 11).Value = 375
$ws.Cells.Item(31, 13).Value = -87
$ws.Cells.Item(109, 8).Value = 1641.3334
$ws.Cells.Item(109, 9).Value = 1641.3334
$ws.Cells.Item(109, 11).Value = 4924.0002
$ws.Cells.Item(109, 13).Value = -3884.0002
$ws.Cells.Item(132, 8).Value = 82385.84
$ws.Cells.Item(132, 10).Value = 3086.5
$ws.Cells.Item(132, 12).Value = 27778.5
$ws.Cells.Item(132, 14).Value = -32838.5
$ws.Cells.Item(135, 8).Value = 728.7143
$ws.Cells.Item(135, 9).Value = 621.4
$ws.Cells.Item(135, 11).Value = 5592.599999999999
$ws.Cells.Item(135, 13).Value = -3057.599999999999
$ws.Cells.Item(141, 8).Value = 7194.2856
$ws.Cells.Item(141, 9).Value = 7194.2856
$ws.Cells.Item(141, 11).Value = 21582.8568
$ws.Cells.Item(141, 13).Value = -16402.8568
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 855.13635
$ws.Cells.Item(102, 9).Value = 684.9474
$ws.Cells.Item(102, 10).Value = 1933
$ws.Cells.Item(102, 11).Value = 684.9474
$ws.Cells.Item(102, 12).Value = 1933
$ws.Cells.Item(102, 13).Value = 937.0526
$ws.Cells.Item(102, 14).Value = -5177
$ws.Cells.Item(122, 8).Value = 1395.5883
$ws.Cells.Item(122, 9).Value = 1599.125
$ws.Cells.Item(122, 10).Value = 1214.6666
$ws.Cells.Item(122, 11).Value = 4797.375
$ws.Cells.Item(122, 12).Value = 3643.9998
$ws.Cells.Item(122, 13).Value = -2347.375
$ws.Cells.Item(122, 14).Value = -8543.9998
$ws.Cells.Item(132, 8).Value = 4520.387
$ws.Cells.Item(132, 9).Value = 2362.4285
$ws.Cells.Item(132, 11).Value = 7087.2855
$ws.Cells.Item(132, 13).Value = -4557.2855
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1709.1428
$ws.Cells.Item(7, 9).Value = 1329.5
$ws.Cells.Item(7, 11).Value = 1329.5
$ws.Cells.Item(7, 13).Value = -1217.5
$ws.Cells.Item(122, 8).Value = 2443.125
$ws.Cells.Item(122, 9).Value = 2440.647
$ws.Cells.Item(122, 10).Value = 2449.1428
$ws.Cells.Item(122, 11).Value = 7321.941
$ws.Cells.Item(122, 12).Value = 7347.428400000001
$ws.Cells.Item(122, 13).Value = -4871.941
$ws.Cells.Item(122, 14).Value = -12247.4284
$ws.Cells.Item(126, 8).Value = 1709.1428
$ws.Cells.Item(126, 9).Value = 1329.5
$ws.Cells.Item(126, 11).Value = 3988.5
$ws.Cells.Item(126, 13).Value = -1518.5
$ws.Cells.Item(132, 8).Value = 4848.5884
$ws.Cells.Item(132, 9).Value = 3062.1667
$ws.Cells.Item(132, 10).Value = 9136
$ws.Cells.Item(132, 11).Value = 9186.500100000001
$ws.Cells.Item(132, 12).Value = 27408
$ws.Cells.Item(132, 13).Value = -6656.500100000001
$ws.Cells.Item(132, 14).Value = -32468
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2973.5
$ws.Cells.Item(96, 9).Value = 1645.5
$ws.Cells.Item(96, 10).Value = 3416.1667
$ws.Cells.Item(96, 11).Value = 1645.5
$ws.Cells.Item(96, 12).Value = 3416.1667
$ws.Cells.Item(96, 13).Value = -272.5
$ws.Cells.Item(96, 14).Value = -6162.1667
$ws.Cells.Item(122, 8).Value = 3076.5
$ws.Cells.Item(122, 9).Value = 3139.9412
$ws.Cells.Item(122, 10).Value = 1998
$ws.Cells.Item(122, 11).Value = 9419.8236
$ws.Cells.Item(122, 12).Value = 5994
$ws.Cells.Item(122, 13).Value = -6969.8236
$ws.Cells.Item(122, 14).Value = -10894
$ws.Cells.Item(126, 8).Value = 5202.3
$ws.Cells.Item(126, 9).Value = 4766.6665
$ws.Cells.Item(126, 11).Value = 14299.9995
$ws.Cells.Item(126, 13).Value = -11829.9995
$ws.Cells.Item(136, 8).Value = 3941.8635
$ws.Cells.Item(136, 9).Value = 1626.8572
$ws.Cells.Item(136, 11).Value = 4880.571599999999
$ws.Cells.Item(136, 13).Value = -2330.571599999999
